$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab_3a_Postulate")

$ws.Range("D4").Value = "Realising globally the right to food"
$ws.Range("D10").Value = "Promoting equal opportunities in society"
$ws.Range("D12").Value = "Reduction the pollution of water with substances"
$ws.Range("D17").Value = "Consolidating public finances – Creating intergenerational equity"
$ws.Range("D18").Value = "Creating favourable investment conditions – Securing long-term prosperity"
$ws.Range("D22").Value = "Shaping the future with new solutions"
$ws.Range("D26").Value = "Guaranteeing mobility – Protecting the environment"
$ws.Range("D30").Value = "Increasing the proportion of sustainable production continuously"
$ws.Range("D33").Value = "Germany's contribution to international climate finance"
$ws.Range("D35").Value = "Conserving species –Protecting habitats"
$ws.Range("D37").Value = "Preventing deforestation and protecting soils world-wide"
